$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2786885245901639
$ws.Range("C2").Value = 0.3934426229508197
$ws.Range("J2").Value = 0.01639344262295082
$ws.Range("P2").Value = 0.2295081967213115
$ws.Range("S2").Value = 0.08196721311475409
$ws.Range("C3").Value = 0.04166666666666666
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.5416666666666666
$ws.Range("S3").Value = 0.375
$ws.Range("B6").Value = 0.02702702702702703
$ws.Range("E6").Value = 0.02702702702702703
$ws.Range("F6").Value = 0.05405405405405406
$ws.Range("J6").Value = 0.2162162162162162
$ws.Range("Q6").Value = 0.1081081081081081
$ws.Range("R6").Value = 0.1081081081081081
$ws.Range("S6").Value = 0.4594594594594595
$ws.Range("B7").Value = 0.1458333333333333
$ws.Range("F7").Value = 0.02083333333333333
$ws.Range("J7").Value = 0.1875
$ws.Range("Q7").Value = 0.0625
$ws.Range("R7").Value = 0.02083333333333333
$ws.Range("S7").Value = 0.5625
$ws.Range("B8").Value = 0.0859375
$ws.Range("F8").Value = 0.046875
$ws.Range("J8").Value = 0.1015625
$ws.Range("O8").Value = 0.0234375
$ws.Range("Q8").Value = 0.09375
$ws.Range("R8").Value = 0.078125
$ws.Range("S8").Value = 0.5703125
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.02272727272727273
$ws.Range("F9").Value = 0.04545454545454546
$ws.Range("J9").Value = 0.04545454545454546
$ws.Range("Q9").Value = 0.06818181818181818
$ws.Range("R9").Value = 0.02272727272727273
$ws.Range("S9").Value = 0.7045454545454546
$ws.Range("B10").Value = 0.09502262443438914
$ws.Range("F10").Value = 0.07692307692307693
$ws.Range("J10").Value = 0.1447963800904978
$ws.Range("O10").Value = 0.01809954751131222
$ws.Range("Q10").Value = 0.1176470588235294
$ws.Range("R10").Value = 0.07239819004524888
$ws.Range("S10").Value = 0.4751131221719457
$ws.Range("G11").Value = 0.2027027027027027
$ws.Range("J11").Value = 0.0945945945945946
$ws.Range("K11").Value = 0.2702702702702703
$ws.Range("L11").Value = 0.3783783783783784
$ws.Range("S11").Value = 0.05405405405405406
$ws.Range("G12").Value = 0.8571428571428571
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("H15").Value = 0.1538461538461539
$ws.Range("I15").Value = 0.1153846153846154
$ws.Range("J15").Value = 0.1923076923076923
$ws.Range("K15").Value = 0.1153846153846154
$ws.Range("M15").Value = 0.03846153846153846
$ws.Range("O15").Value = 0.03846153846153846
$ws.Range("S15").Value = 0.3461538461538461
$ws.Range("H16").Value = 0.3333333333333333
$ws.Range("I16").Value = 0.03703703703703703
$ws.Range("J16").Value = 0.3703703703703703
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("O16").Value = 0.03703703703703703
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.02083333333333333
$ws.Range("H17").Value = 0.1666666666666667
$ws.Range("I17").Value = 0.125
$ws.Range("J17").Value = 0.4166666666666667
$ws.Range("K17").Value = 0.08333333333333333
$ws.Range("H18").Value = 0.125
$ws.Range("I18").Value = 0.09375
$ws.Range("J18").Value = 0.4375
$ws.Range("K18").Value = 0.1875
$ws.Range("M18").Value = 0.03125
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01183431952662722
$ws.Range("H19").Value = 0.3076923076923077
$ws.Range("I19").Value = 0.09171597633136094
$ws.Range("J19").Value = 0.2928994082840237
$ws.Range("K19").Value = 0.1094674556213018
$ws.Range("M19").Value = 0.02366863905325444
$ws.Range("N19").Value = 0.002958579881656805
$ws.Range("O19").Value = 0.03254437869822485
$ws.Range("S19").Value = 0.1272189349112426
